$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cells = @(
    "C4", "D4", "E4", "F4",
    "E5", "F5",
    "B6", "C6", "D6", "E6",
    "D7", "E7",
    "B8", "C8", "D8", "E8",
    "B9", "C9", "D9", "E9",
    "C10", "D10", "E10", "F10",
    "C11", "D11",
    "B12", "C12", "D12", "E12",
    "B13", "C13", "D13", "E13"
)

foreach ($cell in $cells) {
    $old = $ws.Range($cell).Value2
    $new = $old -replace "\.", ","
    $ws.Range($cell).Value2 = $new
}
